$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Extensions list rewrite (paragraph 3)
# ------------------------------------------------------------------
$oldExtList = "Abuse filter, APNG, ImageMap, WikiEditor, BootStrap, CSS, Data, Nuke, MediaWikiChat, InputBox, and ArticleComments.  With these, I should be "
$newExtList = "Abuse filter, ImageMap, WikiEditor, Nuke, YouTube, ConfirmEdit, Spam/TitleBlacklist, CiteThisPage, ParserFunctions, TwoColConflict, Gadgets, InputBox, and possibly more.  With these, I should be "
$d.Content.Find.Execute($oldExtList, $false, $false, $false, $false, $false, $true, 1, $false, $newExtList, 2)

Write-Host "done"
